# "zve test correction #4"
# Fix the ZVE (taxable-income) surcharge formula in column Z: the MIN(...)
# term should be added AFTER the age-factor multiplication, not included
# inside it. i.e.
#   (0.6+(0.02*(T-2005)))*((12*M)+MIN(12*(P+N+0.96*O),2800))
# becomes
#   ((0.6+(0.02*(T-2005)))*(12*M))+MIN(12*(P+N+0.96*O),2800)
# Also updates the view so the selection / scroll position match what the
# workbook was left at after the edit.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Re-enter the corrected formula for every row of the Z column (2-25); each
# cell keeps its own row-relative references, exactly like typing the
# formula once and filling it down over the whole column.
for ($r = 2; $r -le 25; $r++) {
    $formula = "=((0.6+(0.02*(T" + $r + "-2005)))*(12*M" + $r + "))+MIN(12*(P" + $r + "+N" + $r + "+0.96*O" + $r + "),2800)"
    $ws.Cells.Item($r, 26).Formula = $formula
}

# Move the view: scroll so column N is leftmost, and select Z5 (matching
# where the editor ended up after making the correction).
$ws.Activate()
$excel.ActiveWindow.ScrollColumn = 14
$excel.ActiveWindow.ScrollRow = 1
$ws.Range("Z5").Select()
